$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stg_uri_pages")

# Grow the table (Table9 / table10.xml) from A1:E2 to A1:E10 so the new
# rows become part of the structured table, not just loose cell data.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E10"))

# Row 2 (first new table row) filled in completely first.
$ws.Cells.Item(2, 1).Value = "stg_uri_pages"
$ws.Cells.Item(2, 2).Value = "uri_page_id"
$ws.Cells.Item(2, 3).Value = "INTEGER"
$ws.Cells.Item(2, 4).Value = $true
$ws.Cells.Item(2, 5).Value = "table: raw_uri_pages, field: raw_uri_id"

# Remaining rows (3-10), filled column by column.
$colA = @("stg_uri_pages", "stg_uri_pages", "stg_uri_pages", "stg_uri_pages", "stg_uri_pages", "stg_uri_pages", "stg_uri_pages", "stg_uri_pages")
$colB = @("coin_id", "examples_pagination_id", "examples_total_pagination", "examples_start_id", "examples_end_id", "examples_max_id", "uri_link", "ts")
$colC = @("INTEGER", "INTEGER", "INTEGER", "INTEGER", "INTEGER", "INTEGER", "VARCHAR", "TIMESTAMP")
$colD = @($false, $false, $false, $false, $false, $false, $false, $false)
$colE = @("raw_uri_pages", "raw_uri_pages", "raw_uri_pages", "raw_uri_pages", "raw_uri_pages", "raw_uri_pages", "raw_uri_pages", "Database generated")

for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item($i + 3, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item($i + 3, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item($i + 3, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item($i + 3, 4).Value = $colD[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item($i + 3, 5).Value = $colE[$i]
}

# Make stg_uri_pages the active/selected sheet, with B7 selected.
$ws.Activate()
$ws.Range("B7").Select()
